$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (Spanish labels -> snake_case English) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case Spanish connector words (de/del/el/la/las/los/y) in place names ---
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B6").Value = "San José De Gracia"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B22").Value = "Amatenango Del Valle"
$ws.Range("B32").Value = "Comitán De Domínguez"
$ws.Range("B45").Value = "Montecristo De Guerrero"
$ws.Range("B49").Value = "Ocozocoautla De Espinosa"
$ws.Range("B55").Value = "Salto De Agua"
$ws.Range("B56").Value = "San Cristóbal De Las Casas"
$ws.Range("B80").Value = "Guadalupe Y Calvo"
$ws.Range("B81").Value = "Hidalgo Del Parral"
$ws.Range("B89").Value = "San Francisco De Borja"
$ws.Range("B90").Value = "San Francisco De Conchos"
$ws.Range("B91").Value = "San Francisco Del Oro"
$ws.Range("B93").Value = "Valle De Zaragoza"
$ws.Range("B109").Value = "San Juan De Sabinas"
$ws.Range("A119").Value = "Ciudad De México"
$ws.Range("B123").Value = "Cuajimalpa De Morelos"
$ws.Range("B137").Value = "Coneto De Comonfort"
$ws.Range("B149").Value = "Nombre De Dios"
$ws.Range("B153").Value = "San Juan Del Río"
$ws.Range("A161").Value = "Estado De México"
$ws.Range("B161").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B164").Value = "Almoloya De Alquisiras"
$ws.Range("B165").Value = "Almoloya De Juárez"
$ws.Range("B175").Value = "Chapa De Mota"
$ws.Range("B177").Value = "Coacalco De Berriozábal"
$ws.Range("B180").Value = "Ecatepec De Morelos"
$ws.Range("B185").Value = "Ixtapan De La Sal"
$ws.Range("B186").Value = "Ixtapan Del Oro"
$ws.Range("B196").Value = "Naucalpan De Juárez"
$ws.Range("B201").Value = "San Antonio La Isla"
$ws.Range("B202").Value = "San Felipe Del Progreso"
$ws.Range("B203").Value = "San Martín De Las Pirámides"
$ws.Range("B205").Value = "Soyaniquilpan De Juárez"
$ws.Range("B213").Value = "Tenango Del Valle"
$ws.Range("B219").Value = "Tlalnepantla De Baz"
$ws.Range("B224").Value = "Valle De Bravo"
$ws.Range("B225").Value = "Valle De Chalco Solidaridad"
$ws.Range("B226").Value = "Villa De Allende"
$ws.Range("B236").Value = "Apaseo El Grande"
$ws.Range("B242").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B245").Value = "Jaral Del Progreso"
$ws.Range("B252").Value = "Purísima Del Rincón"
$ws.Range("B257").Value = "San Francisco Del Rincón"
$ws.Range("B258").Value = "San Luis De La Paz"
$ws.Range("B259").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B261").Value = "Silao De La Victoria"
$ws.Range("B265").Value = "Valle De Santiago"
$ws.Range("B271").Value = "Acapulco De Juárez"
$ws.Range("B272").Value = "Ajuchitlán Del Progreso"
$ws.Range("B276").Value = "Atenango Del Río"
$ws.Range("B278").Value = "Atoyac De Álvarez"
$ws.Range("B279").Value = "Ayutla De Los Libres"
$ws.Range("B281").Value = "Buenavista De Cuéllar"
$ws.Range("B282").Value = "Chilapa De Álvarez"
$ws.Range("B283").Value = "Chilpancingo De Los Bravo"
$ws.Range("B284").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B289").Value = "Coyuca De Benítez"
$ws.Range("B290").Value = "Coyuca De Catalán"
$ws.Range("B292").Value = "Cuetzala Del Progreso"
$ws.Range("B293").Value = "Cutzamala De Pinzón"
$ws.Range("B297").Value = "Huitzuco De Los Figueroa"
$ws.Range("B298").Value = "Iguala De La Independencia"
$ws.Range("B300").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B301").Value = "Zihuatanejo De Azueta"
$ws.Range("B303").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B305").Value = "Mártir De Cuilapan"
$ws.Range("B315").Value = "Taxco De Alarcón"
$ws.Range("B317").Value = "Técpan De Galeana"
$ws.Range("B319").Value = "Tepecoacuilco De Trujano"
$ws.Range("B320").Value = "Tixtla De Guerrero"
$ws.Range("B323").Value = "Tlapa De Comonfort"
$ws.Range("B332").Value = "Cuautepec De Hinojosa"
$ws.Range("B335").Value = "Huasca De Ocampo"
$ws.Range("B338").Value = "Huejutla De Reyes"
$ws.Range("B341").Value = "Jacala De Ledezma"
$ws.Range("B347").Value = "Omitlán De Juárez"
$ws.Range("B348").Value = "Pachuca De Soto"
$ws.Range("B349").Value = "Progreso De Obregón"
$ws.Range("B352").Value = "Santiago De Anaya"
$ws.Range("B354").Value = "Tenango De Doria"
$ws.Range("B356").Value = "Tepehuacán De Guerrero"
$ws.Range("B357").Value = "Tezontepec De Aldama"
$ws.Range("B361").Value = "Tula De Allende"
$ws.Range("B362").Value = "Tulancingo De Bravo"
$ws.Range("B363").Value = "Zacualtipán De Ángeles"
$ws.Range("B364").Value = "Zapotlán De Juárez"
$ws.Range("B370").Value = "Atemajac De Brizuela"
$ws.Range("B371").Value = "Atotonilco El Alto"
$ws.Range("B372").Value = "Autlán De Navarro"
$ws.Range("B385").Value = "Encarnación De Díaz"
$ws.Range("B390").Value = "Huejuquilla El Alto"
$ws.Range("B391").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B392").Value = "Ixtlahuacán Del Río"
$ws.Range("B395").Value = "Jilotlán De Los Dolores"
$ws.Range("B398").Value = "La Manzanilla De La Paz"
$ws.Range("B399").Value = "Lagos De Moreno"
$ws.Range("B405").Value = "San Diego De Alejandría"
$ws.Range("B406").Value = "San Juan De Los Lagos"
$ws.Range("B407").Value = "San Juanito De Escobedo"
$ws.Range("B410").Value = "San Martín De Bolaños"
$ws.Range("B411").Value = "San Sebastián Del Oeste"
$ws.Range("B412").Value = "Santa María De Los Ángeles"
$ws.Range("B413").Value = "Santa María Del Oro"
$ws.Range("B416").Value = "Tamazula De Gordiano"
$ws.Range("B421").Value = "Teocuitatlán De Corona"
$ws.Range("B422").Value = "Tepatitlán De Morelos"
$ws.Range("B423").Value = "Tizapán El Alto"
$ws.Range("B428").Value = "Unión De San Antonio"
$ws.Range("B429").Value = "Unión De Tula"
$ws.Range("B430").Value = "Valle De Juárez"
$ws.Range("B435").Value = "Yahualica De González Gallo"
$ws.Range("B436").Value = "Zacoalco De Torres"
$ws.Range("B439").Value = "Zapotlán Del Rey"
$ws.Range("B440").Value = "Zapotlán El Grande"
$ws.Range("B499").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B528").Value = "Puente De Ixtla"
$ws.Range("B531").Value = "Tetela Del Volcán"
$ws.Range("B532").Value = "Tlaltizapán De Zapata"
$ws.Range("B538").Value = "Zacualpan De Amilpas"
$ws.Range("B543").Value = "Ixtlán Del Río"
$ws.Range("B556").Value = "Montemorelos"
$ws.Range("B558").Value = "San Nicolás De Los Garza"
$ws.Range("B561").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B564").Value = "Ayoquezco De Aldama"
$ws.Range("B566").Value = "Chiquihuitlán De Benito Juárez"
$ws.Range("B569").Value = "Cuilápam De Guerrero"
$ws.Range("B570").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B571").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B572").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B573").Value = "Ixtlán De Juárez"
$ws.Range("B574").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B578").Value = "Mártires De Tacubaya"
$ws.Range("B580").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B582").Value = "Nejapa De Madero"
$ws.Range("B583").Value = "Oaxaca De Juárez"
$ws.Range("B584").Value = "Ocotlán De Morelos"
$ws.Range("B585").Value = "Putla Villa De Guerrero"
$ws.Range("B597").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B604").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B655").Value = "San Miguel Del Puerto"
$ws.Range("B665").Value = "San Pablo Villa De Mitla"
$ws.Range("B666").Value = "San Pedro El Alto"
$ws.Range("B677").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B678").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B695").Value = "Santa Inés Del Monte"
$ws.Range("B717").Value = "Santiago Del Río"
$ws.Range("B735").Value = "Santo Domingo De Morelos"
$ws.Range("B744").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B745").Value = "Tanetze De Zaragoza"
$ws.Range("B747").Value = "Tataltepec De Valdés"
$ws.Range("B748").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B749").Value = "Teotitlán De Flores Magón"
$ws.Range("B750").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B751").Value = "Tlacolula De Matamoros"
$ws.Range("B752").Value = "Tlalixtac De Cabrera"
$ws.Range("B753").Value = "Totontepec Villa De Morelos"
$ws.Range("B755").Value = "Villa De Chilapa De Díaz"
$ws.Range("B756").Value = "Villa De Etla"
$ws.Range("B757").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B759").Value = "Villa Sola De Vega"
$ws.Range("B760").Value = "Villa Talea De Castro"
$ws.Range("B762").Value = "Zimatlán De Álvarez"
$ws.Range("B773").Value = "Ayotoxco De Guerrero"
$ws.Range("B775").Value = "Chalchicomula De Sesma"
$ws.Range("B793").Value = "Huehuetlán El Chico"
$ws.Range("B796").Value = "Izúcar De Matamoros"
$ws.Range("B802").Value = "Los Reyes De Juárez"
$ws.Range("B819").Value = "San Salvador El Seco"
$ws.Range("B820").Value = "San Salvador El Verde"
$ws.Range("B824").Value = "Tecali De Herrera"
$ws.Range("B830").Value = "Tepanco De López"
$ws.Range("B835").Value = "Tepexi De Rodríguez"
$ws.Range("B840").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B851").Value = "Xayacatlán De Bravo"
$ws.Range("B862").Value = "Amealco De Bonfil"
$ws.Range("B864").Value = "Cadereyta De Montes"
$ws.Range("B867").Value = "Jalpan De Serra"
$ws.Range("B870").Value = "San Juan Del Río"
$ws.Range("B877").Value = "Axtla De Terrazas"
$ws.Range("B880").Value = "Ciudad Del Maíz"
$ws.Range("B885").Value = "Mexquitic De Carmona"
$ws.Range("B889").Value = "San Ciro De Acosta"
$ws.Range("B892").Value = "Santa María Del Río"
$ws.Range("B896").Value = "Villa De Guadalupe"
$ws.Range("B897").Value = "Villa De La Paz"
$ws.Range("B898").Value = "Villa De Ramos"
$ws.Range("B899").Value = "Villa De Reyes"
$ws.Range("B922").Value = "Nacozari De García"
$ws.Range("B949").Value = "Soto La Marina"
$ws.Range("B955").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B959").Value = "Contla De Juan Cuamatzi"
$ws.Range("B962").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B964").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B966").Value = "Papalotla De Xicohténcatl"
$ws.Range("B980").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B984").Value = "Amatlán De Los Reyes"
$ws.Range("B992").Value = "Boca Del Río"
$ws.Range("B994").Value = "Camarón De Tejeda"
$ws.Range("B1009").Value = "Cosamaloapan De Carpio"
$ws.Range("B1021").Value = "Hueyapan De Ocampo"
$ws.Range("B1022").Value = "Ignacio De La Llave"
$ws.Range("B1025").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1026").Value = "Ixhuatlán Del Café"
$ws.Range("B1034").Value = "Juchique De Ferrer"
$ws.Range("B1037").Value = "Landero Y Coss"
$ws.Range("B1039").Value = "Lerdo De Tejada"
$ws.Range("B1042").Value = "Martínez De La Torre"
$ws.Range("B1046").Value = "Mixtla De Altamirano"
$ws.Range("B1060").Value = "Paso De Ovejas"
$ws.Range("B1061").Value = "Paso Del Macho"
$ws.Range("B1064").Value = "Poza Rica De Hidalgo"
$ws.Range("B1072").Value = "Sayula De Alemán"
$ws.Range("B1076").Value = "Tatahuicapan De Juárez"
$ws.Range("B1102").Value = "Vega De Alatorre"
$ws.Range("B1128").Value = "Jiménez Del Teul"
$ws.Range("B1132").Value = "Mezquital Del Oro"
$ws.Range("B1135").Value = "Moyahua De Estrada"
$ws.Range("B1136").Value = "Nochistlán De Mejía"
$ws.Range("B1137").Value = "Noria De Ángeles"
$ws.Range("B1146").Value = "Teúl De González Ortega"
$ws.Range("B1147").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1148").Value = "Trinidad García De La Cadena"

# --- Fix 1-ULP floating point rounding on recomputed percentages ---
$ws.Range("D3").Value = 0.0009532888465204956
$ws.Range("D66").Value = 0.0009532888465204956
$ws.Range("D114").Value = 0.0009532888465204956
$ws.Range("D141").Value = 0.0009532888465204956
$ws.Range("D222").Value = 0.0009532888465204956
$ws.Range("D254").Value = 0.0009532888465204956
$ws.Range("D275").Value = 0.0009532888465204956
$ws.Range("D312").Value = 0.0009532888465204956
$ws.Range("D332").Value = 0.0009532888465204956
$ws.Range("D369").Value = 0.0009532888465204956
$ws.Range("D394").Value = 0.0009532888465204956
$ws.Range("D475").Value = 0.0009532888465204956
$ws.Range("D525").Value = 0.0009532888465204956
$ws.Range("D568").Value = 0.0009532888465204956
$ws.Range("D575").Value = 0.0009532888465204956
$ws.Range("D642").Value = 0.0009532888465204956
$ws.Range("D747").Value = 0.0009532888465204956
$ws.Range("D792").Value = 0.0009532888465204956
$ws.Range("D840").Value = 0.0009532888465204956
$ws.Range("D908").Value = 0.0009532888465204956
$ws.Range("D961").Value = 0.0009532888465204956
$ws.Range("D1029").Value = 0.0009532888465204956
$ws.Range("D1072").Value = 0.0009532888465204956
$ws.Range("D1151").Value = 0.0009532888465204956
$ws.Range("D125").Value = 0.009941440827999456
$ws.Range("D770").Value = 0.009396704344273456
$ws.Range("D1138").Value = 0.009396704344273456

# --- Remove trailing footnote rows (now-removed metadata block) ---
$ws.Range("A1157:A1161").EntireRow.Delete()

Write-Output "done"